$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Forces a numeric-looking string to be stored as text (not coerced to a
    # number) the way real Excel requires, then strips the resulting
    # "Text" number-format style back off so the cell's style index is left
    # untouched (matches original formatting / avoids spurious style diffs).
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.915.52"
$ws.Range("E2").Value = "  +1.70%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.645.46"
$ws.Range("E3").Value = "  +1.87%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "213.55"
$ws.Range("E5").Value = "  +1.16%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.26%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.07%  "

# Row 8 - Solana
Set-TextValue "D8" "23.39"
$ws.Range("E8").Value = "  +2.70%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.39%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.54%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0871"
$ws.Range("E11").Value = "  -1.61%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.879.15"
$ws.Range("E12").Value = "  +1.85%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.645.50"
$ws.Range("E13").Value = "  +1.86%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.564"
$ws.Range("E15").Value = "  +3.03%  "

# Row 16 - Litecoin
Set-TextValue "D16" "65.57"
$ws.Range("E16").Value = "  +0.73%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.934.71"

# Row 18 - BitcoinCash
Set-TextValue "D18" "231.31"
$ws.Range("E18").Value = "  -0.76%  "

# Row 19 - was Chainlink, now ShibaInu
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  +0.99%  "

# Row 20 - was ShibaInu, now Chainlink
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D20" "7.66"
$ws.Range("E20").Value = "  +2.01%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.02%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  +4.41%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +1.97%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +3.66%  "

# Row 25 - Monero
Set-TextValue "D25" "152.37"
$ws.Range("E25").Value = "  +1.43%  "

# Row 26 - Cosmos
Set-TextValue "D26" "6.91"
$ws.Range("E26").Value = "  +0.89%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "15.74"
$ws.Range("E28").Value = "  +1.46%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.07%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.68%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +0.59%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.33"
$ws.Range("E32").Value = "  +2.06%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.442.68"
$ws.Range("E33").Value = "  -1.72%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +0.26%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +1.64%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.14%  "

# Row 37 - ARBITRUM
$ws.Range("E37").Value = "  +3.29%  "

# Row 38 - was VeChain, now TrustWalletToken
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D38" "0.934"
$ws.Range("E38").Value = "  -3.47%  "

# Row 39 - was TrustWalletToken, now VeChain
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.0169"
$ws.Range("E39").Value = "  +1.18%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "  +0.53%  "

# Row 41 - Aave
Set-TextValue "D41" "69.11"
$ws.Range("E41").Value = "  +3.27%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  +3.68%  "

# Row 43 - PaxDollar
$ws.Range("E43").Value = "  -0.05%  "

# Row 44 - mCoin
$ws.Range("E44").Value = "  -0.02%  "

# Row 45 - RenderToken
$ws.Range("E45").Value = "  +6.08%  "

# Row 46 - FraxShare
$ws.Range("E46").Value = "  +3.50%  "

# Row 47 - MXToken
Set-TextValue "D47" "2.21"
$ws.Range("E47").Value = "  +0.70%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "1.787.78"
$ws.Range("E48").Value = "  +1.54%  "

# Row 49 - Quant
Set-TextValue "D49" "89.07"
$ws.Range("E49").Value = "  +2.64%  "

# Row 50 - BabyDogeCoin
$ws.Range("D50").Value = "0.0₆0105"
$ws.Range("E50").Value = "  +0.10%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  +0.32%  "
